# Update "paises.xlsx" (COVID-19 country stats) with refreshed source data:
#   - New snapshot timestamp
#   - Refreshed numbers for a few countries (Estados Unidos, Noruega)
#   - Germany ("Alemania") jumps above France ("Francia") in the ranking,
#     carrying fresh numbers while France keeps its previous row's figures
#   - Ghana moves up in the ranking (now just after Taiwan, before Jordania)
#     with fresh numbers, pushing Jordania..Mauricio down by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 00:52"

# Row updates: Country + [Casos totales, Nuevos casos, Casos activos,
#              Recuperados, Casos criticos, Muertes hoy, Muertes]
$updates = @(
    @{ Row = 4;   Pais = "Estados Unidos"; Vals = @(465240, 30313, 25139, 423589, 9959, 1724, 16512) },
    @{ Row = 7;   Pais = "Alemania";       Vals = @(118181, 4885, 52407, 63167, 4895, 258, 2607) },
    @{ Row = 8;   Pais = "Francia";        Vals = @(117749, 4799, 23206, 82333, 7066, 1341, 12210) },
    @{ Row = 26;  Pais = "Noruega";        Vals = @(6211, 169, 32, 6071, 78, 7, 108) },
    @{ Row = 97;  Pais = "Ghana";          Vals = @(378, 65, 34, 338, 2, 0, 6) },
    @{ Row = 98;  Pais = "Jordania";       Vals = @(372, 14, 161, 204, 5, 1, 7) },
    @{ Row = 99;  Pais = "Reunion";        Vals = @(362, 0, 40, 322, 4, 0, 0) },
    @{ Row = 100; Pais = "Honduras";       Vals = @(343, 31, 6, 314, 10, 1, 23) },
    @{ Row = 101; Pais = "Malta";          Vals = @(337, 38, 16, 319, 4, 1, 2) },
    @{ Row = 102; Pais = "San Marino";     Vals = @(333, 25, 49, 250, 14, 0, 34) },
    @{ Row = 103; Pais = "Banglades";      Vals = @(330, 112, 33, 276, 1, 1, 21) },
    @{ Row = 104; Pais = "Mauricio";       Vals = @(314, 41, 23, 284, 3, 0, 7) }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Pais
    for ($i = 0; $i -lt $u.Vals.Count; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $u.Vals[$i]
    }
}
